$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"  = -12.672
    "C6"  = -12.239
    "C7"  = -12.89
    "C16" = -13.093
    "C20" = -12.13
    "C28" = -12.766
    "C29" = -12.129
    "C32" = -12.608
    "C40" = -12.237
    "C46" = -13.654
    "C51" = -11.296
    "C52" = -11.382
    "C57" = -13.818
    "C59" = -12.442
    "C62" = -13.211
    "C66" = -11.713
    "C73" = -12.199
    "C74" = -11.982
    "C92" = -10.753
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
